# Historical data update: insert 9 new trading-day rows (2019-11-18 .. 2019-11-28)
# immediately before the existing 2019-11-29 row, shifting all subsequent rows
# down by 9 (old row 124 -> new row 133, ... old row 196 -> new row 205).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 124:196 down by 9 rows to make room for the new data.
$ws.Rows("124:132").Insert()

# New rows to insert, in order: timestamp, date, id, name, open, high, low, close, vol
$newRows = @(
    @(1574035200, "2019-11-18", "6633", "LHI", 0.92,              0.93,              0.92,   0.93,              6406000),
    @(1574121600, "2019-11-19", "6633", "LHI", 0.925,             0.9350000000000001,0.92,   0.925,             4317200),
    @(1574208000, "2019-11-20", "6633", "LHI", 0.92,              0.9350000000000001,0.915,  0.925,             2883600),
    @(1574294400, "2019-11-21", "6633", "LHI", 0.93,              0.9350000000000001,0.925,  0.93,              3444400),
    @(1574380800, "2019-11-22", "6633", "LHI", 0.93,              0.9399999999999999,0.92,   0.93,              4833400),
    @(1574640000, "2019-11-25", "6633", "LHI", 0.93,              0.9399999999999999,0.925,  0.9350000000000001,8432900),
    @(1574726400, "2019-11-26", "6633", "LHI", 0.9350000000000001,0.9350000000000001,0.9,    0.9,               10425300),
    @(1574812800, "2019-11-27", "6633", "LHI", 0.905,             0.905,             0.855,  0.86,              9301100),
    @(1574899200, "2019-11-28", "6633", "LHI", 0.86,              0.865,             0.85,   0.85,              3693300)
)

$startRow = 124
$endRow = 132

# Columns B (date) and C (id) must stay plain text, exactly like every other
# row in the sheet, instead of being auto-converted to a date serial / number.
$ws.Range("B$($startRow):C$($endRow)").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}

# Restore the default (unstyled) cell style on the two text columns now that
# the values are safely stored as text, so no stray number-format style is
# left behind on these cells.
$ws.Range("B$($startRow):C$($endRow)").Style = "Normal"
